# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holding detail) right before
#    the "总计" (totals) sheet.
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (row 2) for the 2022-Q1 quarter, shifting the existing rows down,
#    and keep the running index column (A) sequential.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" detail sheet
# ---------------------------------------------------------------------
# NOTE: worksheet variables resolve by live position, not stable object
# identity - once Worksheets.Add() shuffles positions, a previously
# captured reference can silently start pointing at a different sheet.
# So: do the Add()+rename first, then re-fetch every handle we need by
# name right before using it.
$totalSheetForInsert = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetForInsert)
$newSheet.Name = "2022-Q1"

$newSheet   = $wb.Worksheets.Item("2022-Q1")
$template   = $wb.Worksheets.Item("2021-Q4")

# Carry over the look (bold/centered/bordered header + index-column
# style) from an existing quarter sheet with the same 8-column layout.
# (Column A has no header cell, so copy B:H for row 1 only.)
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:H2").Copy()
$newSheet.Range("A2:H11").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @("000751", "嘉实新兴产业股票",               "81.04", "90.66", "2.72", "2.2043", 10),
    @("004450", "嘉实前沿科技沪港深股票",           "21.72", "88.32", "6.39", "1.3879", 3),
    @("004616", "中欧电子信息产业沪港深股票A",       "14.54", "92.26", "5.77", "0.8390", 3),
    @("005763", "中欧电子信息产业沪港深股票C",       "7.73",  "92.26", "5.77", "0.4460", 3),
    @("012447", "华夏互联网龙头混合型证券投资基金A", "3.32",  "83.95", "6.55", "0.2175", 7),
    @("010016", "华夏科技前沿6个月定期开放混合A",     "7.78",  "91.92", "2.38", "0.1852", 8),
    @("001759", "嘉实成长增强灵活配置混合",           "4.59",  "90.80", "3.12", "0.1432", 10),
    @("012448", "华夏互联网龙头混合型证券投资基金C", "1.28",  "83.95", "6.55", "0.0838", 7),
    @("010017", "华夏科技前沿6个月定期开放混合C",     "2.39",  "91.92", "2.38", "0.0569", 8),
    @("013903", "国泰君安信息行业混合",               "0.25",  "84.06", "3.71", "0.0093", 3)
)

$r = 2
foreach ($fund in $fundRows) {
    $newSheet.Range("A$r").Value = ($r - 2)

    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $fund[0]

    $newSheet.Range("C$r").Value = $fund[1]

    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value = $fund[2]

    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value = $fund[3]

    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value = $fund[4]

    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value = $fund[5]

    $newSheet.Range("H$r").Value = $fund[6]

    $r++
}

# ---------------------------------------------------------------------
# 2. Prepend a 2022-Q1 summary row to the "总计" sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# The inserted row copies formatting down from row 1 (header), whose
# A-column cell is blank/unstyled - so A2 loses the bold/centred index
# style every other row in this column carries. Restore it from A3.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 10
$totalSheet.Range("D2").Value = 5.57

# Renumber the running index in column A for the rows pushed down.
for ($row = 3; $row -le 7; $row++) {
    $totalSheet.Range("A$row").Value = $row - 2
}
